# Red Alarm (Virtual Boy) - Plan.xlsx
# "Updated set to better align with writing policy"
#
# The achievement-description helper formulas in column I of the
# "Achievements" sheet wrapped the difficulty/reset clause in literal
# parentheses, e.g.:
#     H2 & " (" & IF(...) & IF(...) & ")"
# -> "Complete stage 1 (on easy+ difficulty)"
#
# The policy update drops the parentheses and just separates the clauses
# with a plain space:
#     H2 & " " & IF(...) & IF(...)
# -> "Complete stage 1 on easy+ difficulty"
#
# Every row 2..61 uses the same per-row template (only the row number
# changes), so rewrite them all in one pass. Everything downstream
# (Checklist!C2:C61, Text!A3:A.. via INDIRECT) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Achievements")

for ($r = 2; $r -le 61; $r++) {
    $formula = 'H' + $r + '&" "&IF(F' + $r + '="Easy","on easy+ difficulty",IF(F' + $r + '="Normal","on normal+ difficulty",IF(F' + $r + '="Hard","on hard difficulty","")))&IF(G' + $r + '="Game Over",", resets on new game",IF(G' + $r + '="Session",", resets on new session",""))'
    $ws.Range("I" + $r).Formula = "=" + $formula
}

# View-state bookkeeping that went along with the edit: the editor had been
# working on column I (selection ends on I47:I61, the last edited block)
# before switching away from "Achievements" and leaving "Leaderboards" as
# the active tab when the file was saved.
$ws.Activate()
$ws.Range("I47:I61").Select()

$wsLeaderboards = $wb.Worksheets.Item("Leaderboards")
$wsLeaderboards.Activate()
